$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 currently holds a half-finished "Test" course record. Complete it:
#  - courseName: "Test" -> "Java"
#  - property:   "optional" -> "compulsory"
#  - credit:     "1.0" -> "3.0"
$ws.Range("C12").Value = "Java"
$ws.Range("D12").Value = "compulsory"

# "3.0" looks numeric, and a bare assignment would make Excel store it as
# the number 3 instead of literal text (this column holds text everywhere
# else, e.g. "4.0", "1.0", "2.5"). Prefix with an apostrophe to force a
# text entry, then lift the style back off the cell (PasteSpecial formats
# from an untouched neighbor) so E12 keeps the sheet's default styling
# instead of picking up an explicit quote-prefixed style.
$ws.Range("E12").Value = "'3.0"
$ws.Range("D12").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$excel.CutCopyMode = 0
